$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 10.31202652642125
$ws.Cells.Item(2, 4).Value = 4.958632881046337
$ws.Cells.Item(2, 5).Value = 11.34184065356047
$ws.Cells.Item(2, 6).Value = 49.94632990694338
$ws.Cells.Item(2, 7).Value = 65.6837915430914
$ws.Cells.Item(2, 8).Value = 23.51358782473558
$ws.Cells.Item(2, 10).Value = 10.04799720429031
$ws.Cells.Item(2, 12).Value = 8.476595989078547
$ws.Cells.Item(2, 13).Value = 62.76036394584125

# Row 3
$ws.Cells.Item(3, 3).Value = 10.46997525035151
$ws.Cells.Item(3, 4).Value = 4.853770489368898
$ws.Cells.Item(3, 5).Value = 11.3047279243419
$ws.Cells.Item(3, 6).Value = 50.66104340379878
$ws.Cells.Item(3, 7).Value = 66.5947747329275
$ws.Cells.Item(3, 8).Value = 23.79444759018024
$ws.Cells.Item(3, 10).Value = 10.12531108345997
$ws.Cells.Item(3, 12).Value = 8.40303260247274
$ws.Cells.Item(3, 13).Value = 59.71620204156128

# Row 4
$ws.Cells.Item(4, 3).Value = 10.57043731653633
$ws.Cells.Item(4, 4).Value = 4.787654474737342
$ws.Cells.Item(4, 5).Value = 11.28517581504225
$ws.Cells.Item(4, 6).Value = 51.12693918830038
$ws.Cells.Item(4, 7).Value = 67.19903524652912
$ws.Cells.Item(4, 8).Value = 23.97705369120794
$ws.Cells.Item(4, 10).Value = 10.1759863206904
$ws.Cells.Item(4, 12).Value = 8.358854466338695
$ws.Cells.Item(4, 13).Value = 57.76072275469672

# Row 5
$ws.Cells.Item(5, 3).Value = 10.61226107045064
$ws.Cells.Item(5, 4).Value = 4.760294021643734
$ws.Cells.Item(5, 5).Value = 11.27801896261407
$ws.Cells.Item(5, 6).Value = 51.32350514396327
$ws.Cells.Item(5, 7).Value = 67.45631633007577
$ws.Cells.Item(5, 8).Value = 24.05399717029794
$ws.Cells.Item(5, 10).Value = 10.19743549727441
$ws.Cells.Item(5, 12).Value = 8.34111157164908
$ws.Cells.Item(5, 13).Value = 56.94271650078443

# Row 6
$ws.Cells.Item(6, 3).Value = 10.61925957172018
$ws.Cells.Item(6, 4).Value = 4.755726143290437
$ws.Cells.Item(6, 5).Value = 11.27687945418722
$ws.Cells.Item(6, 6).Value = 51.35654737682586
$ws.Cells.Item(6, 7).Value = 67.49969728492896
$ws.Cells.Item(6, 8).Value = 24.06692570929252
$ws.Cells.Item(6, 10).Value = 10.20104515303457
$ws.Cells.Item(6, 12).Value = 8.338181404521681
$ws.Cells.Item(6, 13).Value = 56.80563005748249

# Row 7
$ws.Cells.Item(7, 3).Value = 10.57099777220725
$ws.Cells.Item(7, 4).Value = 4.787287149378913
$ws.Cells.Item(7, 5).Value = 11.28507601627535
$ws.Cells.Item(7, 6).Value = 51.12956309840501
$ws.Cells.Item(7, 7).Value = 67.2024606538227
$ws.Cells.Item(7, 8).Value = 23.97808116452255
$ws.Cells.Item(7, 10).Value = 10.17627236723652
$ws.Cells.Item(7, 12).Value = 8.358614111986922
$ws.Cells.Item(7, 13).Value = 57.74977555035024

# Row 8
$ws.Cells.Item(8, 3).Value = 10.36577050845062
$ws.Cells.Item(8, 4).Value = 4.922845520705034
$ws.Cells.Item(8, 5).Value = 11.32837116556951
$ws.Cells.Item(8, 6).Value = 50.18708822026562
$ws.Cells.Item(8, 7).Value = 65.98842001941748
$ws.Cells.Item(8, 8).Value = 23.60830522117704
$ws.Cells.Item(8, 10).Value = 10.07398554837281
$ws.Cells.Item(8, 12).Value = 8.451032537351507
$ws.Cells.Item(8, 13).Value = 61.72897388584272

# Row 9
$ws.Cells.Item(9, 3).Value = 9.990498632698323
$ws.Cells.Item(9, 4).Value = 5.174239710691023
$ws.Cells.Item(9, 5).Value = 11.43905517487627
$ws.Cells.Item(9, 6).Value = 48.55817081545181
$ws.Cells.Item(9, 7).Value = 63.9761253971764
$ws.Cells.Item(9, 8).Value = 22.96492407925648
$ws.Cells.Item(9, 10).Value = 9.899172622581021
$ws.Cells.Item(9, 12).Value = 8.639631583979194
$ws.Cells.Item(9, 13).Value = 68.82682144465983

# Row 10
$ws.Cells.Item(10, 3).Value = 9.73071387404811
$ws.Cells.Item(10, 4).Value = 5.349273162108151
$ws.Cells.Item(10, 5).Value = 11.53623400763577
$ws.Cells.Item(10, 6).Value = 47.50193745502528
$ws.Cells.Item(10, 7).Value = 62.73999238665871
$ws.Cells.Item(10, 8).Value = 22.54379105091039
$ws.Cells.Item(10, 10).Value = 9.786971664743639
$ws.Cells.Item(10, 12).Value = 8.782034884683092
$ws.Cells.Item(10, 13).Value = 73.59096680787377

# Row 11
$ws.Cells.Item(11, 3).Value = 9.615844013188216
$ws.Cells.Item(11, 4).Value = 5.426641186095941
$ws.Cells.Item(11, 5).Value = 11.58391119507982
$ws.Cells.Item(11, 6).Value = 47.05362523673597
$ws.Cells.Item(11, 7).Value = 62.23442969791233
$ws.Cells.Item(11, 8).Value = 22.36382324394052
$ws.Cells.Item(11, 10).Value = 9.739582138413114
$ws.Cells.Item(11, 12).Value = 8.847514567297281
$ws.Cells.Item(11, 13).Value = 75.65759815347219

# Row 12
$ws.Cells.Item(12, 3).Value = 9.572809290217078
$ws.Cells.Item(12, 4).Value = 5.455602094803395
$ws.Cells.Item(12, 5).Value = 11.60246555801229
$ws.Cells.Item(12, 6).Value = 46.88864358717667
$ws.Cells.Item(12, 7).Value = 62.05152050191265
$ws.Cells.Item(12, 8).Value = 22.29738253093575
$ws.Cells.Item(12, 10).Value = 9.722173933267197
$ws.Cells.Item(12, 12).Value = 8.872399211797379
$ws.Cells.Item(12, 13).Value = 76.42554965698798

# Row 13
$ws.Cells.Item(13, 3).Value = 9.582057140950601
$ws.Cells.Item(13, 4).Value = 5.449380003391684
$ws.Cells.Item(13, 5).Value = 11.59844728954132
$ws.Cells.Item(13, 6).Value = 46.92395989324684
$ws.Cells.Item(13, 7).Value = 62.090527430881
$ws.Cells.Item(13, 8).Value = 22.31161503882356
$ws.Cells.Item(13, 10).Value = 9.725899007003646
$ws.Cells.Item(13, 12).Value = 8.867036087455043
$ws.Cells.Item(13, 13).Value = 76.26081072754431

# Row 14
$ws.Cells.Item(14, 3).Value = 9.612294289588359
$ws.Cells.Item(14, 4).Value = 5.429030631728786
$ws.Cells.Item(14, 5).Value = 11.58542764674766
$ws.Cells.Item(14, 6).Value = 47.03995525897893
$ws.Cells.Item(14, 7).Value = 62.21920834906278
$ws.Cells.Item(14, 8).Value = 22.35832263127016
$ws.Cells.Item(14, 10).Value = 9.738139116185284
$ws.Cells.Item(14, 12).Value = 8.849560126185265
$ws.Cells.Item(14, 13).Value = 75.72107213031212

# Row 15
$ws.Cells.Item(15, 3).Value = 9.630875490320383
$ws.Cells.Item(15, 4).Value = 5.416521861194052
$ws.Cells.Item(15, 5).Value = 11.57751788823147
$ws.Cells.Item(15, 6).Value = 47.11163375564801
$ws.Cells.Item(15, 7).Value = 62.29915228043674
$ws.Cells.Item(15, 8).Value = 22.38715620962582
$ws.Cells.Item(15, 10).Value = 9.745706869544867
$ws.Cells.Item(15, 12).Value = 8.838866831221047
$ws.Cells.Item(15, 13).Value = 75.38855629256084

# Row 16
$ws.Cells.Item(16, 3).Value = 9.738286462251827
$ws.Cells.Item(16, 4).Value = 5.34417051172666
$ws.Cells.Item(16, 5).Value = 11.5331881128736
$ws.Cells.Item(16, 6).Value = 47.53189633569338
$ws.Cells.Item(16, 7).Value = 62.77420579536124
$ws.Cells.Item(16, 8).Value = 22.55578929504595
$ws.Cells.Item(16, 10).Value = 9.790143169857586
$ws.Cells.Item(16, 12).Value = 8.777768706457861
$ws.Cells.Item(16, 13).Value = 73.4538674763146

# Row 17
$ws.Cells.Item(17, 3).Value = 9.805018621493064
$ws.Cells.Item(17, 4).Value = 5.299198461095814
$ws.Cells.Item(17, 5).Value = 11.50688281696615
$ws.Cells.Item(17, 6).Value = 47.79806259201351
$ws.Cells.Item(17, 7).Value = 63.08043332870344
$ws.Cells.Item(17, 8).Value = 22.6622407159403
$ws.Cells.Item(17, 10).Value = 9.818346908902829
$ws.Cells.Item(17, 12).Value = 8.740457786082491
$ws.Cells.Item(17, 13).Value = 72.24109413698821

# Row 18
$ws.Cells.Item(18, 3).Value = 9.843713235750494
$ws.Cells.Item(18, 4).Value = 5.273120087220806
$ws.Cells.Item(18, 5).Value = 11.49207927673161
$ws.Cells.Item(18, 6).Value = 47.95417250665287
$ws.Cells.Item(18, 7).Value = 63.26188621153148
$ws.Cells.Item(18, 8).Value = 22.7245587486563
$ws.Cells.Item(18, 10).Value = 9.834911843854742
$ws.Cells.Item(18, 12).Value = 8.719064077606776
$ws.Cells.Item(18, 13).Value = 71.53407648053549

# Row 19
$ws.Cells.Item(19, 3).Value = 9.856868507603014
$ws.Cells.Item(19, 4).Value = 5.264254442946844
$ws.Cells.Item(19, 5).Value = 11.48712310530837
$ws.Cells.Item(19, 6).Value = 48.00754240377895
$ws.Cells.Item(19, 7).Value = 63.32422549465014
$ws.Cells.Item(19, 8).Value = 22.74584460334407
$ws.Cells.Item(19, 10).Value = 9.840578983655018
$ws.Cells.Item(19, 12).Value = 8.711832325870219
$ws.Cells.Item(19, 13).Value = 71.29307329287627

# Row 20
$ws.Cells.Item(20, 3).Value = 9.797882656926514
$ws.Cells.Item(20, 4).Value = 5.304007809508649
$ws.Cells.Item(20, 5).Value = 11.50964925351816
$ws.Cells.Item(20, 6).Value = 47.76941527232626
$ws.Cells.Item(20, 7).Value = 63.04728175210915
$ws.Cells.Item(20, 8).Value = 22.65079569088467
$ws.Cells.Item(20, 10).Value = 9.815309000038143
$ws.Cells.Item(20, 12).Value = 8.744422795648745
$ws.Cells.Item(20, 13).Value = 72.37117665852844

# Row 21
$ws.Cells.Item(21, 3).Value = 9.603400407925486
$ws.Cells.Item(21, 4).Value = 5.435016965657497
$ws.Cells.Item(21, 5).Value = 11.58923825282049
$ws.Cells.Item(21, 6).Value = 47.0057534358418
$ws.Cells.Item(21, 7).Value = 62.18117688249806
$ws.Cells.Item(21, 8).Value = 22.34455675713032
$ws.Cells.Item(21, 10).Value = 9.734529218951415
$ws.Cells.Item(21, 12).Value = 8.854690918524415
$ws.Cells.Item(21, 13).Value = 75.88000479316996

# Row 22
$ws.Cells.Item(22, 3).Value = 9.478993857245756
$ws.Cells.Item(22, 4).Value = 5.518671648431772
$ws.Cells.Item(22, 5).Value = 11.64416928244677
$ws.Cells.Item(22, 6).Value = 46.53462651774783
$ws.Cells.Item(22, 7).Value = 61.66507239327527
$ws.Cells.Item(22, 8).Value = 22.15439530232657
$ws.Cells.Item(22, 10).Value = 9.684872385318887
$ws.Cells.Item(22, 12).Value = 8.927270530701135
$ws.Cells.Item(22, 13).Value = 78.08786628756377

# Row 23
$ws.Cells.Item(23, 3).Value = 9.545148983844213
$ws.Cells.Item(23, 4).Value = 5.474207368817307
$ws.Cells.Item(23, 5).Value = 11.61458453329714
$ws.Cells.Item(23, 6).Value = 46.78346093234111
$ws.Cells.Item(23, 7).Value = 61.93582807898827
$ws.Cells.Item(23, 8).Value = 22.25496049453872
$ws.Cells.Item(23, 10).Value = 9.711083835791619
$ws.Cells.Item(23, 12).Value = 8.888490275511797
$ws.Cells.Item(23, 13).Value = 76.91734327033039

# Row 24
$ws.Cells.Item(24, 3).Value = 9.801107799603956
$ws.Cells.Item(24, 4).Value = 5.301834197806478
$ws.Cells.Item(24, 5).Value = 11.50839755176331
$ws.Cells.Item(24, 6).Value = 47.78235712082086
$ws.Cells.Item(24, 7).Value = 63.06225279256117
$ws.Cells.Item(24, 8).Value = 22.65596650833618
$ws.Cells.Item(24, 10).Value = 9.816681348655194
$ws.Cells.Item(24, 12).Value = 8.742630036899094
$ws.Cells.Item(24, 13).Value = 72.31239679370468

# Row 25
$ws.Cells.Item(25, 3).Value = 10.08917491754509
$ws.Cells.Item(25, 4).Value = 5.107870033168248
$ws.Cells.Item(25, 5).Value = 11.40633599277107
$ws.Cells.Item(25, 6).Value = 48.97468926642887
$ws.Cells.Item(25, 7).Value = 64.47930105339498
$ws.Cells.Item(25, 8).Value = 23.13005491158086
$ws.Cells.Item(25, 10).Value = 9.943653959021132
$ws.Cells.Item(25, 12).Value = 8.476595989078547
$ws.Cells.Item(25, 13).Value = 66.98461720017285
